$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-08-06 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-07 Monday", 2) | Out-Null
$d.Content.Find.Execute("11÷9=1, 2", $true, $false, $false, $false, $false, $true, 1, $false, "59÷2=29, 1", 2) | Out-Null
$d.Content.Find.Execute("96÷6=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "59÷2=29, 1", 2) | Out-Null
$d.Content.Find.Execute("23÷6=3, 5", $true, $false, $false, $false, $false, $true, 1, $false, "60÷2=30, 0", 2) | Out-Null
$d.Content.Find.Execute("19÷4=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "31÷5=6, 1", 2) | Out-Null
$d.Content.Find.Execute("41÷4=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "62÷2=31, 0", 2) | Out-Null
$d.Content.Find.Execute("88÷2=44, 0", $true, $false, $false, $false, $false, $true, 1, $false, "79÷7=11, 2", 2) | Out-Null
$d.Content.Find.Execute("39÷9=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "34÷5=6, 4", 2) | Out-Null
$d.Content.Find.Execute("95÷8=11, 7", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=6, 0", 2) | Out-Null
$d.Content.Find.Execute("38÷3=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "74÷3=24, 2", 2) | Out-Null
$d.Content.Find.Execute("42÷6=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "45÷8=5, 5", 2) | Out-Null
$d.Content.Find.Execute("72÷5=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "99÷6=16, 3", 2) | Out-Null
$d.Content.Find.Execute("89÷5=17, 4", $true, $false, $false, $false, $false, $true, 1, $false, "31÷8=3, 7", 2) | Out-Null
$d.Content.Find.Execute("11÷4=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "23÷4=5, 3", 2) | Out-Null
$d.Content.Find.Execute("68÷2=34, 0", $true, $false, $false, $false, $false, $true, 1, $false, "22÷6=3, 4", 2) | Out-Null
$d.Content.Find.Execute("85÷9=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "87÷4=21, 3", 2) | Out-Null
$d.Content.Find.Execute("39÷3=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "13÷4=3, 1", 2) | Out-Null
$d.Content.Find.Execute("71÷8=8, 7", $true, $false, $false, $false, $false, $true, 1, $false, "69÷9=7, 6", 2) | Out-Null
$d.Content.Find.Execute("89÷6=14, 5", $true, $false, $false, $false, $false, $true, 1, $false, "10÷7=1, 3", 2) | Out-Null
$d.Content.Find.Execute("39÷5=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "18÷2=9, 0", 2) | Out-Null
$d.Content.Find.Execute("34÷9=3, 7", $true, $false, $false, $false, $false, $true, 1, $false, "66÷3=22, 0", 2) | Out-Null
$d.Content.Find.Execute("50÷8=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "83÷6=13, 5", 2) | Out-Null
$d.Content.Find.Execute("65÷2=32, 1", $true, $false, $false, $false, $false, $true, 1, $false, "94÷7=13, 3", 2) | Out-Null
$d.Content.Find.Execute("54÷3=18, 0", $true, $false, $false, $false, $false, $true, 1, $false, "74÷3=24, 2", 2) | Out-Null
$d.Content.Find.Execute("67÷6=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=5, 0", 2) | Out-Null
